$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new K1 header cell the same base look (bold/white on dark fill)
# as the other "white text" headers (G1/H1/I1/J1) before we touch its fill,
# so the later accent-fill tweak reuses the existing bold-white font run
# instead of inventing a duplicate one.
$ws.Range("I1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Header row (row 1): drop the "* " required-field prefix ---
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Employee No."
$ws.Range("C1").Value = "First Name"
$ws.Range("D1").Value = "Last Name"
$ws.Range("E1").Value = "Middle Name"
$ws.Range("F1").Value = "Email"
$ws.Range("G1").Value = "Password"
$ws.Range("H1").Value = "Access Keys"
$ws.Range("I1").Value = "Role"

# New "Location(s)" column inserted at J, "User Reset" shifts from J to K
$ws.Range("J1").Value = "Location(s)"
$ws.Range("K1").Value = "User Reset"

# --- Sample data row (row 2) ---
$ws.Range("J2").Value = "METRO MANILA, CENTRAL LUZON"
$ws.Range("K2").Value = 1

# --- New note row (row 3) ---
$ws.Range("A3").Value = "* delete this line and above sample data on actual uploading"

# --- Formatting: new header fill (blue accent) for Middle Name / Location(s) / User Reset ---
$ws.Range("E1").Interior.ThemeColor = 5
$ws.Range("J1").Interior.ThemeColor = 5
$ws.Range("K1").Interior.ThemeColor = 5

# --- Formatting: red bold warning text for the note row ---
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Font.Color = 255

# --- Widen the Location(s) column to fit its new contents ---
$ws.Columns.Item(10).ColumnWidth = 28.43

# --- Selection moves to the new note row ---
$ws.Range("A3").Select() | Out-Null

Write-Output "edit applied"
